# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Cerberus Profits workbook
# (commit: "chore: update Sheets via scheduled runner")

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 306.9
$ws.Range("I4").Value = 230
$ws.Range("K4").Value = 230
$ws.Range("M4").Value = -116
$ws.Range("H18").Value = 1529.25
$ws.Range("I18").Value = 1529.25
$ws.Range("K18").Value = 1529.25
$ws.Range("M18").Value = -1245.25
$ws.Range("H92").Value = 985.8
$ws.Range("I92").Value = 699.0714
$ws.Range("K92").Value = 699.0714
$ws.Range("M92").Value = 548.9286
$ws.Range("H107").Value = 661.7273
$ws.Range("J107").Value = 376.875
$ws.Range("L107").Value = 376.875
$ws.Range("N107").Value = -4216.875
$ws.Range("H138").Value = 3095
$ws.Range("J138").Value = 2954.0967
$ws.Range("L138").Value = 8862.2901
$ws.Range("N138").Value = -19142.2901
$ws.Range("H141").Value = 4937
$ws.Range("I141").Value = 4096.6665
$ws.Range("K141").Value = 12289.9995
$ws.Range("M141").Value = -7109.999500000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1295.9333
$ws.Range("I32").Value = 798.7143
$ws.Range("K32").Value = 798.7143
$ws.Range("M32").Value = -511.7143
$ws.Range("H61").Value = 7964.05
$ws.Range("I61").Value = 6402.294
$ws.Range("K61").Value = 6402.294
$ws.Range("M61").Value = -6190.294
$ws.Range("H74").Value = 1683.6
$ws.Range("I74").Value = 1175.8
$ws.Range("K74").Value = 1175.8
$ws.Range("M74").Value = -301.8
$ws.Range("H77").Value = 1683.6
$ws.Range("I77").Value = 1175.8
$ws.Range("K77").Value = 5879
$ws.Range("M77").Value = -1511
$ws.Range("H122").Value = 4756.6665
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1319.9524
$ws.Range("I132").Value = 802.8627300000001
$ws.Range("K132").Value = 2408.58819
$ws.Range("M132").Value = 121.4118099999996
$ws.Range("H136").Value = 7964.05
$ws.Range("I136").Value = 6402.294
$ws.Range("K136").Value = 19206.882
$ws.Range("M136").Value = -16656.882

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2424.44
$ws.Range("I105").Value = 2365.1904
$ws.Range("J105").Value = 2735.5
$ws.Range("K105").Value = 2365.1904
$ws.Range("L105").Value = 2735.5
$ws.Range("M105").Value = -618.1904
$ws.Range("N105").Value = -6229.5
$ws.Range("H134").Value = 7698.15
$ws.Range("I134").Value = 6620.654
$ws.Range("K134").Value = 19861.962
$ws.Range("M134").Value = -17326.962

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 39202.5
$ws.Range("J28").Value = 39202.5
$ws.Range("L28").Value = 39202.5
$ws.Range("N28").Value = -39692.5
$ws.Range("H134").Value = 4560.9736
$ws.Range("I134").Value = 3916.25
$ws.Range("J134").Value = 10041.125
$ws.Range("K134").Value = 11748.75
$ws.Range("L134").Value = 30123.375
$ws.Range("M134").Value = -9213.75
$ws.Range("N134").Value = -35193.375
$ws.Range("H140").Value = 114209.6
$ws.Range("J140").Value = 114209.6
$ws.Range("L140").Value = 114209.6
$ws.Range("N140").Value = -124569.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 1562.6154
$ws.Range("J40").Value = 2490.5
$ws.Range("L40").Value = 9962
$ws.Range("N40").Value = -10100
$ws.Range("H92").Value = 8152.25
$ws.Range("J92").Value = 10003
$ws.Range("L92").Value = 30009
$ws.Range("N92").Value = -32505
$ws.Range("H122").Value = 1878.4286
$ws.Range("I122").Value = 354.66666
$ws.Range("K122").Value = 3191.99994
$ws.Range("M122").Value = -741.9999399999997
$ws.Range("H129").Value = 15155212
$ws.Range("J129").Value = 18522838
$ws.Range("L129").Value = 55568514
$ws.Range("N129").Value = -55578514

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 32499.75
$ws.Range("J15").Value = 32499.75
$ws.Range("L15").Value = 32499.75
$ws.Range("N15").Value = -33075.75
$ws.Range("H81").Value = 32499.75
$ws.Range("J81").Value = 32499.75
$ws.Range("L81").Value = 32499.75
$ws.Range("N81").Value = -34495.75
$ws.Range("H84").Value = 32499.75
$ws.Range("J84").Value = 32499.75
$ws.Range("L84").Value = 97499.25
$ws.Range("N84").Value = -107483.25
$ws.Range("H97").Value = 529.2857
$ws.Range("I97").Value = 582.63635
$ws.Range("K97").Value = 582.63635
$ws.Range("M97").Value = -86.63634999999999
$ws.Range("H102").Value = 3767.8064
$ws.Range("J102").Value = 3914.3333
$ws.Range("L102").Value = 3914.3333
$ws.Range("N102").Value = -7158.3333
$ws.Range("H104").Value = 44999.668
$ws.Range("J104").Value = 44999.668
$ws.Range("L104").Value = 44999.668
$ws.Range("N104").Value = -51987.668
$ws.Range("H113").Value = 3866.1428
$ws.Range("I113").Value = 1842
$ws.Range("J113").Value = 4675.8
$ws.Range("K113").Value = 1842
$ws.Range("L113").Value = 4675.8
$ws.Range("M113").Value = 328
$ws.Range("N113").Value = -9015.799999999999
$ws.Range("H132").Value = 1179.3256
$ws.Range("I132").Value = 1031.4849
$ws.Range("K132").Value = 3094.4547
$ws.Range("M132").Value = -564.4546999999998
$ws.Range("H134").Value = 18000
$ws.Range("J134").Value = 18000
$ws.Range("L134").Value = 54000
$ws.Range("N134").Value = -59070

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1886.6666
$ws.Range("I16").Value = 1968.95
$ws.Range("K16").Value = 1968.95
$ws.Range("M16").Value = -1798.95
$ws.Range("H30").Value = 1375
$ws.Range("I30").Value = 1375
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1375
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1267
$ws.Range("N30").ClearContents()
$ws.Range("H68").Value = 2545.5625
$ws.Range("I68").Value = 2414.889
$ws.Range("J68").Value = 2713.5715
$ws.Range("K68").Value = 2414.889
$ws.Range("L68").Value = 2713.5715
$ws.Range("M68").Value = -1665.889
$ws.Range("N68").Value = -4211.5715
$ws.Range("H71").Value = 2545.5625
$ws.Range("I71").Value = 2414.889
$ws.Range("J71").Value = 2713.5715
$ws.Range("K71").Value = 12074.445
$ws.Range("L71").Value = 13567.8575
$ws.Range("M71").Value = -8330.445
$ws.Range("N71").Value = -21055.8575
$ws.Range("H105").Value = 21110.5
$ws.Range("J105").Value = 21110.5
$ws.Range("L105").Value = 21110.5
$ws.Range("N105").Value = -28098.5
$ws.Range("H106").Value = 8040
$ws.Range("J106").Value = 8040
$ws.Range("L106").Value = 8040
$ws.Range("N106").Value = -10564
$ws.Range("H122").Value = 4986.5
$ws.Range("I122").Value = 2598
$ws.Range("J122").Value = 7375
$ws.Range("K122").Value = 7794
$ws.Range("L122").Value = 22125
$ws.Range("M122").Value = -5344
$ws.Range("N122").Value = -27025
$ws.Range("H135").Value = 96124.125
$ws.Range("I135").Value = 100000
$ws.Range("J135").Value = 95570.42999999999
$ws.Range("K135").Value = 100000
$ws.Range("L135").Value = 95570.42999999999
$ws.Range("M135").Value = -94930
$ws.Range("N135").Value = -105710.43
$ws.Range("H136").Value = 2125.675
$ws.Range("I136").Value = 1307.0952
$ws.Range("K136").Value = 3921.2856
$ws.Range("M136").Value = -1371.2856

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9333.333000000001
$ws.Range("J81").Value = 4499.75
$ws.Range("L81").Value = 8999.5
$ws.Range("N81").Value = -11121.5
$ws.Range("H84").Value = 9333.333000000001
$ws.Range("J84").Value = 4499.75
$ws.Range("L84").Value = 44997.5
$ws.Range("N84").Value = -55605.5
$ws.Range("H96").Value = 3118.2307
$ws.Range("I96").Value = 2833.6667
$ws.Range("K96").Value = 2833.6667
$ws.Range("M96").Value = -1460.6667
$ws.Range("H105").Value = 7500
$ws.Range("J105").Value = 7500
$ws.Range("L105").Value = 7500
$ws.Range("N105").Value = -14488
$ws.Range("H107").Value = 893.0357
$ws.Range("I107").Value = 834.9048
$ws.Range("K107").Value = 2504.7144
$ws.Range("M107").Value = -584.7143999999998
$ws.Range("H140").Value = 321327.88
$ws.Range("J140").Value = 321327.88
$ws.Range("L140").Value = 321327.88
$ws.Range("N140").Value = -331687.88
$ws.Range("H141").Value = 81760.53999999999
$ws.Range("J141").Value = 81760.53999999999
$ws.Range("L141").Value = 81760.53999999999
$ws.Range("N141").Value = -92120.53999999999
